{"js": "// Append the M2Doc \"version mismatch\" warning marker to the document body's\n// (only) paragraph, right after the existing (empty) run \u2014 without merging\n// into it, matching the 4 distinct new <w:r> runs from the diff:\n//   1. \"    \" (plain, 4 spaces)\n//   2. \"<---\" (orange, size 32 half-pt == 16pt, lightGray highlight)\n//   3. \"M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0\" (same formatting as #2)\n//   4. \"    \" (plain, 4 spaces)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Collapsed range at the very end of the paragraph (after its last run).\nconst endRange = paragraph.getRange(\"End\");\n\n// Using insertOoxml (Range.InsertXML under the hood) inserts the fragment as\n// literal, already-built runs instead of \"typing\" text \u2014 which keeps each\n// <w:r> distinct (no merge with the pre-existing empty run or with each\n// other), exactly like the target OOXML diff.\nconst warningRunsOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">    </w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"FFA500\"/><w:sz w:val=\"32\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>&lt;---</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"FFA500\"/><w:sz w:val=\"32\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">    </w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nendRange.insertOoxml(warningRunsOoxml, \"End\");\nawait context.sync();\n", "ps1": "# Append the M2Doc \"version mismatch\" warning marker after the existing\n# (empty) run of the document's last paragraph, matching the 4 distinct new\n# <w:r> runs added by the diff:\n#   1. \"    \" (plain, 4 spaces)\n#   2. \"<---\" (orange FFA500, size 32 half-pt == 16pt, lightGray highlight)\n#   3. \"M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0\" (same formatting as #2)\n#   4. \"    \" (plain, 4 spaces)\n#\n# Word's COM \"typing\" APIs (Range.InsertAfter/InsertBefore, Selection.TypeText,\n# ...) coalesce a newly inserted run into an adjacent run that already carries\n# identical (here: empty/default) run properties \u2014 which would silently fold\n# the brand-new \"    \" run back into the pre-existing empty run instead of\n# keeping them as two separate <w:r> elements. Range.InsertXML does not do\n# that (it splices in literal XML), but it always inserts/replaces whole\n# paragraphs, so we replace the *entire* current (only) paragraph \u2014 mark\n# included \u2014 with a paragraph holding the old empty run plus the 4 new runs,\n# then delete the stray extra paragraph mark that operation leaves behind.\n\n$d = $word.ActiveDocument\n\n$paragraphs = $d.Paragraphs\n$paragraphCountBefore = $paragraphs.Count\n$targetParagraph = $paragraphs.Item($paragraphCountBefore)\n$targetRange = $targetParagraph.Range\n\n$warningRunsOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:t/></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">    </w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"FFA500\"/><w:sz w:val=\"32\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>&lt;---</w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"FFA500\"/><w:sz w:val=\"32\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">    </w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n# Replace the whole paragraph (incl. its end-of-paragraph mark) with the\n# reconstructed one; InsertXML on a Range spanning a full paragraph mark\n# inserts our new paragraph and leaves one extra (empty) paragraph behind.\n$targetRange.InsertXML($warningRunsOoxml)\n\n$paragraphsAfter = $d.Paragraphs\nif ($paragraphsAfter.Count -gt $paragraphCountBefore) {\n    $rebuiltParagraph = $paragraphsAfter.Item($paragraphCountBefore)\n    $rebuiltEnd = $rebuiltParagraph.Range.End\n    $lastParagraph = $paragraphsAfter.Item($paragraphsAfter.Count)\n    $strayMark = $d.Range($rebuiltEnd - 1, $lastParagraph.Range.End - 1)\n    $strayMark.Delete()\n}\n"}
